# dialogue, task 추가 중
# - Remove the unused empty "Sheet2" / "Sheet3" tabs
# - Add two new sample/dialogue rows (7 & 8) to the "a4" sheet

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the two empty placeholder sheets
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# Work on the "a4" sheet (first tab) and keep it the active/selected one
$ws = $wb.Worksheets.Item("a4")
$ws.Activate()

# New row 7: asd / 1 / asd / asdasd / asd
$ws.Range("A7").Value = "asd"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "asd"
$ws.Range("D7").Value = "asdasd"
$ws.Range("E7").Value = "asd"

# New row 8: aaa / 2 / (blank) / aaa
$ws.Range("A8").Value = "aaa"
$ws.Range("B8").Value = 2
$ws.Range("D8").Value = "aaa"

# Leave the selection on D8, matching the final cell edited
$ws.Range("D8").Select()
